$wb = $excel.ActiveWorkbook

$wsSPN = $wb.Worksheets.Item("SPN")
$wsITI = $wb.Worksheets.Item("ITI")

# Column C ("Semana") values change from the text "Semana 05" to the plain number 5
$wsSPN.Range("C2:C26").Value = 5
$wsITI.Range("C2:C56").Value = 5

$wb.Save()
